$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates taken from the commit diff (price/volume refresh, plus a
# 2-row reorder of Binance-PegBSC-USD / WrappedeETH at rows 27-28).
# D (Price) / E (Volume 1h) are forced to Text ("@") before the value is
# written so numeric-looking strings (e.g. "25.30", "0.0000170") keep their
# exact formatting instead of being auto-coerced into numbers.
$updates = @(
    @{ Cell = 'D2'; Value = '61.791.95' },
    @{ Cell = 'E2'; Value = '  -0.85%  ' },
    @{ Cell = 'D3'; Value = '2.398.75' },
    @{ Cell = 'E3'; Value = '  -1.13%  ' },
    @{ Cell = 'E4'; Value = '  +0.02%  ' },
    @{ Cell = 'D5'; Value = '560.17' },
    @{ Cell = 'E5'; Value = '  +0.62%  ' },
    @{ Cell = 'D6'; Value = '141.64' },
    @{ Cell = 'E6'; Value = '  -1.30%  ' },
    @{ Cell = 'E7'; Value = '  +0.12%  ' },
    @{ Cell = 'D8'; Value = '0.531' },
    @{ Cell = 'E8'; Value = '  -0.43%  ' },
    @{ Cell = 'E9'; Value = '  -1.52%  ' },
    @{ Cell = 'E10'; Value = '  -1.95%  ' },
    @{ Cell = 'D11'; Value = '5.21' },
    @{ Cell = 'E11'; Value = '  -3.58%  ' },
    @{ Cell = 'D12'; Value = '0.347' },
    @{ Cell = 'E12'; Value = '  -1.36%  ' },
    @{ Cell = 'D13'; Value = '25.30' },
    @{ Cell = 'E13'; Value = '  -3.95%  ' },
    @{ Cell = 'D14'; Value = '0.0000170' },
    @{ Cell = 'E14'; Value = '  -2.21%  ' },
    @{ Cell = 'D15'; Value = '2.832.05' },
    @{ Cell = 'E15'; Value = '  -1.10%  ' },
    @{ Cell = 'D16'; Value = '61.798.15' },
    @{ Cell = 'E16'; Value = '  -0.63%  ' },
    @{ Cell = 'D17'; Value = '2.397.33' },
    @{ Cell = 'E17'; Value = '  -1.17%  ' },
    @{ Cell = 'D18'; Value = '11.13' },
    @{ Cell = 'E18'; Value = '  +0.09%  ' },
    @{ Cell = 'D19'; Value = '320.10' },
    @{ Cell = 'E19'; Value = '  -1.33%  ' },
    @{ Cell = 'D20'; Value = '6.78' },
    @{ Cell = 'E20'; Value = '  +0.45%  ' },
    @{ Cell = 'D21'; Value = '4.09' },
    @{ Cell = 'E21'; Value = '  -1.91%  ' },
    @{ Cell = 'D22'; Value = '0.999' },
    @{ Cell = 'D23'; Value = '65.41' },
    @{ Cell = 'E23'; Value = '  +0.66%  ' },
    @{ Cell = 'E24'; Value = '  -4.75%  ' },
    @{ Cell = 'D25'; Value = '8.64' },
    @{ Cell = 'D26'; Value = '560.34' },
    @{ Cell = 'E26'; Value = '  -1.59%  ' },
    @{ Cell = 'B27'; Value = 'WrappedeETH' },
    @{ Cell = 'C27'; Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth' },
    @{ Cell = 'D27'; Value = '2.519.41' },
    @{ Cell = 'E27'; Value = '  -0.96%  ' },
    @{ Cell = 'B28'; Value = 'Binance-PegBSC-USD' },
    @{ Cell = 'C28'; Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd' },
    @{ Cell = 'D28'; Value = '0.992' },
    @{ Cell = 'E28'; Value = '  -0.85%  ' },
    @{ Cell = 'D29'; Value = '0.0₃0921' },
    @{ Cell = 'E29'; Value = '  -2.38%  ' },
    @{ Cell = 'D30'; Value = '8.11' },
    @{ Cell = 'E30'; Value = '  -3.57%  ' },
    @{ Cell = 'D31'; Value = '1.37' },
    @{ Cell = 'E31'; Value = '  -5.80%  ' },
    @{ Cell = 'E32'; Value = '  -1.69%  ' },
    @{ Cell = 'D33'; Value = '1.85' },
    @{ Cell = 'E33'; Value = '  -0.45%  ' },
    @{ Cell = 'E34'; Value = '  -4.94%  ' },
    @{ Cell = 'E35'; Value = '  +0.14%  ' },
    @{ Cell = 'D36'; Value = '4.73' },
    @{ Cell = 'E36'; Value = '  -2.10%  ' },
    @{ Cell = 'D37'; Value = '152.01' },
    @{ Cell = 'E37'; Value = '  +1.50%  ' },
    @{ Cell = 'D38'; Value = '5.39' },
    @{ Cell = 'E38'; Value = '  -5.96%  ' },
    @{ Cell = 'E39'; Value = '  -1.96%  ' },
    @{ Cell = 'D40'; Value = '18.43' },
    @{ Cell = 'E40'; Value = '  -2.03%  ' },
    @{ Cell = 'D41'; Value = '1.77' },
    @{ Cell = 'E41'; Value = '  -6.27%  ' },
    @{ Cell = 'E42'; Value = '  -0.07%  ' },
    @{ Cell = 'D43'; Value = '146.89' },
    @{ Cell = 'E43'; Value = '  -3.15%  ' },
    @{ Cell = 'D44'; Value = '2.21' },
    @{ Cell = 'E44'; Value = '  -5.37%  ' },
    @{ Cell = 'E45'; Value = '  -1.80%  ' },
    @{ Cell = 'D46'; Value = '0.0526' },
    @{ Cell = 'E46'; Value = '  -3.32%  ' },
    @{ Cell = 'D47'; Value = '19.71' },
    @{ Cell = 'E47'; Value = '  -3.52%  ' },
    @{ Cell = 'D48'; Value = '0.585' },
    @{ Cell = 'E48'; Value = '  -0.81%  ' },
    @{ Cell = 'D49'; Value = '0.0915' },
    @{ Cell = 'E49'; Value = '  -0.24%  ' },
    @{ Cell = 'D50'; Value = '0.0223' },
    @{ Cell = 'E50'; Value = '  -2.33%  ' },
    @{ Cell = 'E51'; Value = '  +0.31%  ' }
)

foreach ($u in $updates) {
    $range = $ws.Range($u.Cell)
    if ($u.Cell -match "^[DE]\d+$") {
        $range.NumberFormat = "@"
    }
    $range.Value = $u.Value
}
